$d = $word.ActiveDocument

$replacements = @(
    @{old = "62÷5=12, 2"; new = "17÷8=2, 1"},
    @{old = "50÷4=12, 2"; new = "13÷9=1, 4"},
    @{old = "97÷4=24, 1"; new = "57÷8=7, 1"},
    @{old = "55÷3=18, 1"; new = "69÷2=34, 1"},
    @{old = "11÷7=1, 4"; new = "11÷4=2, 3"},
    @{old = "62÷7=8, 6"; new = "31÷3=10, 1"},
    @{old = "53÷5=10, 3"; new = "67÷8=8, 3"},
    @{old = "76÷4=19, 0"; new = "89÷9=9, 8"},
    @{old = "55÷9=6, 1"; new = "68÷2=34, 0"},
    @{old = "12÷7=1, 5"; new = "17÷2=8, 1"},
    @{old = "48÷7=6, 6"; new = "98÷2=49, 0"},
    @{old = "56÷2=28, 0"; new = "20÷6=3, 2"},
    @{old = "62÷8=7, 6"; new = "64÷5=12, 4"},
    @{old = "22÷7=3, 1"; new = "50÷3=16, 2"},
    @{old = "73÷7=10, 3"; new = "21÷6=3, 3"},
    @{old = "19÷2=9, 1"; new = "45÷8=5, 5"},
    @{old = "87÷5=17, 2"; new = "85÷2=42, 1"},
    @{old = "41÷7=5, 6"; new = "10÷2=5, 0"},
    @{old = "42÷3=14, 0"; new = "70÷3=23, 1"},
    @{old = "11÷3=3, 2"; new = "55÷4=13, 3"},
    @{old = "12÷8=1, 4"; new = "88÷9=9, 7"},
    @{old = "76÷6=12, 4"; new = "97÷7=13, 6"},
    @{old = "68÷3=22, 2"; new = "57÷7=8, 1"},
    @{old = "75÷9=8, 3"; new = "22÷5=4, 2"},
    @{old = "62÷9=6, 8"; new = "14÷8=1, 6"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
